# Updates to figs and tables list
# Adds two new rows (23 and 24) to the "List of all figures and tables"
# sheet describing the ctenophore-first / sponge-first species tree
# supplementary figures, matching the style used by the existing
# "Online Supplementary Material" rows (21-22) but with the pink/red
# fill used elsewhere in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from the formatting already used by row 22 (Online Supplementary
# Material, pink-ish fill) and stamp it onto the two new rows so the
# borders / fonts / wrap text all line up, then tweak the fill colour and
# vertical alignment to match the new style.
$ws.Range("A22:E22").Copy() | Out-Null
$ws.Range("A23:E24").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Row 23: ctenophore-first species tree
$ws.Range("A23").Value = "Supplementary Figure X"
$ws.Range("B23").Value = "Online Supplementary Material"
$ws.Range("C23").Value = "Species tree with ctenophore first."
$ws.Range("D23").Value = "To be prepared"

# Row 24: sponge-first species tree
$ws.Range("A24").Value = "Supplementary Figure X"
$ws.Range("B24").Value = "Online Supplementary Material"
$ws.Range("C24").Value = "Species tree with sponge first."
$ws.Range("D24").Value = "To be prepared"

# New rows use a pink fill (RGB FF9999) instead of the salmon fill copied
# from row 22, and are not vertically centred (default/bottom alignment).
$ws.Range("A23:E24").Interior.Color = 10066431
$ws.Range("A23:E24").VerticalAlignment = -4107  # xlBottom (default)

# Match the cursor position recorded in the saved workbook.
$ws.Range("B29").Select() | Out-Null
